# Update the "Generate Report for Handback" timestamps.
# These cells store plain text timestamps (t="s" shared strings) even though
# their cell style applies a date/time display format, so we must write them
# as text to avoid Excel re-interpreting the string as a serial date value.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-08-21 10:51:29"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-08-21 10:51:25"
$zhcn.Range("K3").Value = "2016-08-21 10:51:43"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("K3").Value = "2016-08-21 10:51:49"
